# Fixed issues with individual counties
#
# 1) Columns A ("month") and B ("year") had their values swapped in the
#    source data (A held the year, B held the month number). Swap the
#    values back so A contains the month number and B contains the year.
# 2) Add four new "grade" columns (grade_total, grade_distance,
#    grade_visitation, grade_encounters) plus five encounter-frequency
#    label columns (NEVER, RARELY, SOMETIMES, FREQUENTLY, ALWAYS) with
#    their header labels and constant values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# --- 1) Swap A/B values for every data row (row 1 is the header row) ---
for ($r = 2; $r -le $lastRow; $r++) {
    $monthCell = $ws.Cells.Item($r, 1)
    $yearCell = $ws.Cells.Item($r, 2)

    $monthValue = $monthCell.Value()
    $yearValue = $yearCell.Value()

    $monthCell.Value = $yearValue
    $yearCell.Value = $monthValue
}

# --- 2) Add the new header columns (H1:P1) ---
$newHeaders = @(
    "grade_total",
    "grade_distance",
    "grade_visitation",
    "grade_encounters",
    "NEVER",
    "RARELY",
    "SOMETIMES",
    "FREQUENTLY",
    "ALWAYS"
)

$firstNewCol = 8  # column H
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $firstNewCol + $i).Value = $newHeaders[$i]
}

# --- 3) Fill the new columns with their (constant) row values ---
$newValues = @(0, 0, 0, 0, 1.06, 1.091, 1.126, 1.227, 1.496)

for ($r = 2; $r -le $lastRow; $r++) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        $ws.Cells.Item($r, $firstNewCol + $i).Value = $newValues[$i]
    }
}
